$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C5 time value increases by 10 minutes: 12:00:00 PM -> 12:10:00 PM
$ws.Range("C5").Value = 0.5 + (10/1440)
